$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append the new row 17 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "Vragen over samenwerking"
$logs.Range("B17").Value = "mailmind.test@zohomail.eu"
$logs.Range("C17").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D17").Value = "Overig"
$logs.Range("F17").Value = "2025-06-18 15:30:11"
$logs.Range("G17").Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row ---
$catFormats = $logs.Range("D2:D16").FormatConditions
for ($i = 1; $i -le $catFormats.Count; $i++) {
    $catFormats.Item($i).ModifyAppliesToRange($logs.Range("D2:D17"))
}

$answeredFormats = $logs.Range("G2:G16").FormatConditions
for ($i = 1; $i -le $answeredFormats.Count; $i++) {
    $answeredFormats.Item($i).ModifyAppliesToRange($logs.Range("G2:G17"))
}

# --- Update the "Dashboard" sheet: bump the "Overig" count from 7 to 8 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 8
